$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 73569.36
$ws.Range("I6").Value = 85807.586
$ws.Range("K6").Value = 257422.758
$ws.Range("M6").Value = -257310.758

$ws.Range("H92").Value = 1242.95
$ws.Range("I92").Value = 715
$ws.Range("J92").Value = 5994.5
$ws.Range("K92").Value = 715
$ws.Range("L92").Value = 5994.5
$ws.Range("M92").Value = 533
$ws.Range("N92").Value = -8490.5

$ws.Range("H99").Value = 73980150
$ws.Range("I99").Value = 2976666
$ws.Range("K99").Value = 8929998
$ws.Range("M99").Value = -8928500

$ws.Range("H101").Value = 2674684.8
$ws.Range("I101").Value = 4546244.5
$ws.Range("J101").Value = 1028.2858
$ws.Range("K101").Value = 13638733.5
$ws.Range("L101").Value = 3084.8574
$ws.Range("M101").Value = -13637111.5
$ws.Range("N101").Value = -6328.857400000001

$ws.Range("H138").Value = 3042.0708
$ws.Range("J138").Value = 3502.1775
$ws.Range("L138").Value = 10506.5325
$ws.Range("N138").Value = -20786.5325

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 20782
$ws.Range("I2").Value = 33970.332
$ws.Range("J2").Value = 999.5
$ws.Range("K2").Value = 33970.332
$ws.Range("L2").Value = 999.5
$ws.Range("M2").Value = -33857.332
$ws.Range("N2").Value = -1225.5

$ws.Range("H32").Value = 11987.833
$ws.Range("I32").Value = 10927.794
$ws.Range("J32").Value = 34248.668
$ws.Range("K32").Value = 10927.794
$ws.Range("L32").Value = 34248.668
$ws.Range("M32").Value = -10640.794
$ws.Range("N32").Value = -34822.668

$ws.Range("H45").Value = 2450.55
$ws.Range("I45").Value = 2206
$ws.Range("J45").Value = 2477.7222
$ws.Range("K45").Value = 2206
$ws.Range("L45").Value = 2477.7222
$ws.Range("M45").Value = -1829
$ws.Range("N45").Value = -3231.7222

$ws.Range("H61").Value = 2238.0625
$ws.Range("I61").Value = 1371.1482
$ws.Range("K61").Value = 1371.1482
$ws.Range("M61").Value = -1159.1482

$ws.Range("H102").Value = 1754.3549
$ws.Range("I102").Value = 1396.72
$ws.Range("K102").Value = 1396.72
$ws.Range("M102").Value = 225.28

$ws.Range("H116").Value = 20782
$ws.Range("I116").Value = 33970.332
$ws.Range("J116").Value = 999.5
$ws.Range("K116").Value = 33970.332
$ws.Range("L116").Value = 999.5
$ws.Range("M116").Value = -31676.332
$ws.Range("N116").Value = -5587.5

$ws.Range("H136").Value = 2238.0625
$ws.Range("I136").Value = 1371.1482
$ws.Range("K136").Value = 4113.444600000001
$ws.Range("M136").Value = -1563.444600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 20782
$ws.Range("I3").Value = 33970.332
$ws.Range("J3").Value = 999.5
$ws.Range("K3").Value = 33970.332
$ws.Range("L3").Value = 999.5
$ws.Range("M3").Value = -33856.332
$ws.Range("N3").Value = -1227.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 73.59999999999999
$ws.Range("I7").Value = 38.25
$ws.Range("J7").Value = 86.454544
$ws.Range("K7").Value = 38.25
$ws.Range("L7").Value = 86.454544
$ws.Range("M7").Value = 74.75
$ws.Range("N7").Value = -312.454544

$ws.Range("H31").Value = 24308.861
$ws.Range("I31").Value = 1887.2727
$ws.Range("J31").Value = 59542.785
$ws.Range("K31").Value = 1887.2727
$ws.Range("L31").Value = 59542.785
$ws.Range("M31").Value = -1592.2727
$ws.Range("N31").Value = -60132.785

$ws.Range("H34").Value = 24308.861
$ws.Range("I34").Value = 1887.2727
$ws.Range("J34").Value = 59542.785
$ws.Range("K34").Value = 1887.2727
$ws.Range("L34").Value = 59542.785
$ws.Range("M34").Value = -1685.2727
$ws.Range("N34").Value = -59946.785

$ws.Range("H58").Value = 3755.4546
$ws.Range("I58").Value = 3674.7222
$ws.Range("K58").Value = 3674.7222
$ws.Range("M58").Value = -3471.7222

$ws.Range("H132").Value = 3962.7273
$ws.Range("I132").Value = 3756.6458
$ws.Range("K132").Value = 11269.9374
$ws.Range("M132").Value = -8739.937399999999

$ws.Range("H134").Value = 3377.347
$ws.Range("I134").Value = 3445.5
$ws.Range("K134").Value = 10336.5
$ws.Range("M134").Value = -7801.5

$ws.Range("H136").Value = 3755.4546
$ws.Range("I136").Value = 3674.7222
$ws.Range("K136").Value = 11024.1666
$ws.Range("M136").Value = -8474.1666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 158.6
$ws.Range("I2").Value = 140.1
$ws.Range("J2").Value = 177.1
$ws.Range("K2").Value = 840.5999999999999
$ws.Range("L2").Value = 1062.6
$ws.Range("M2").Value = -727.5999999999999
$ws.Range("N2").Value = -1288.6

$ws.Range("H59").Value = 1452.5
$ws.Range("I59").Value = 1452.5
$ws.Range("K59").Value = 4357.5
$ws.Range("M59").Value = -3817.5

$ws.Range("H131").Value = 21740902
$ws.Range("J131").Value = 1962.4103
$ws.Range("L131").Value = 5887.2309
$ws.Range("N131").Value = -15967.2309

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2079.366
$ws.Range("I102").Value = 2110.3784
$ws.Range("J102").Value = 1792.5
$ws.Range("K102").Value = 2110.3784
$ws.Range("L102").Value = 1792.5
$ws.Range("M102").Value = -488.3784000000001
$ws.Range("N102").Value = -5036.5

$ws.Range("H126").Value = 7365.769
$ws.Range("I126").Value = 5293.68
$ws.Range("J126").Value = 11065.929
$ws.Range("K126").Value = 15881.04
$ws.Range("L126").Value = 33197.787
$ws.Range("M126").Value = -13411.04
$ws.Range("N126").Value = -38137.787

$ws.Range("H132").Value = 25743
$ws.Range("I132").Value = 24972.596
$ws.Range("K132").Value = 74917.788
$ws.Range("M132").Value = -72387.788

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 5999
$ws.Range("J3").Value = 5999
$ws.Range("L3").Value = 5999
$ws.Range("N3").Value = -6223

$ws.Range("H15").Value = 5999
$ws.Range("J15").Value = 5999
$ws.Range("L15").Value = 5999
$ws.Range("N15").Value = -6339

$ws.Range("H42").Value = 36961.5
$ws.Range("I42").Value = 36961.5
$ws.Range("K42").Value = 36961.5
$ws.Range("M42").Value = -36398.5

$ws.Range("H43").Value = 1000000
$ws.Range("I43").Value = 1000000
$ws.Range("K43").Value = 1000000
$ws.Range("M43").Value = -999807

$ws.Range("H49").Value = 36961.5
$ws.Range("I49").Value = 36961.5
$ws.Range("K49").Value = 36961.5
$ws.Range("M49").Value = -36814.5

$ws.Range("H100").Value = 3420.4285
$ws.Range("I100").Value = 3095.6
$ws.Range("J100").Value = 4232.5
$ws.Range("K100").Value = 3095.6
$ws.Range("L100").Value = 4232.5
$ws.Range("M100").Value = -2554.6
$ws.Range("N100").Value = -5314.5

$ws.Range("H136").Value = 50541.06
$ws.Range("I136").Value = 3391.3635
$ws.Range("J136").Value = 136982.17
$ws.Range("K136").Value = 10174.0905
$ws.Range("L136").Value = 410946.51
$ws.Range("M136").Value = -7624.0905
$ws.Range("N136").Value = -416046.51

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 619.10345
$ws.Range("I100").Value = 528.2273
$ws.Range("J100").Value = 904.7143
$ws.Range("K100").Value = 1056.4546
$ws.Range("L100").Value = 1809.4286
$ws.Range("M100").Value = -515.4546
$ws.Range("N100").Value = -2891.4286

$ws.Range("H132").Value = 1648.5253
$ws.Range("I132").Value = 945.3585
$ws.Range("J132").Value = 2458.6956
$ws.Range("K132").Value = 2836.0755
$ws.Range("L132").Value = 7376.0868
$ws.Range("M132").Value = -306.0754999999999
$ws.Range("N132").Value = -12436.0868

$ws.Range("H136").Value = 2929.2656
$ws.Range("I136").Value = 1561.6923
$ws.Range("K136").Value = 4685.0769
$ws.Range("M136").Value = -2135.0769
